# Applies the cryptos.xlsx data refresh described in the commit diff.
# D-column "Price" values are numeric-looking strings that Excel would
# otherwise auto-convert to floating point numbers (losing exact text and
# introducing binary-float rounding noise), so we force text storage via
# NumberFormat "@" before the assignment, then restore the default "Normal"
# style afterwards so no stray style index is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.854.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.218.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  -1.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.554"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("E9").Value = "  -5.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0776"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.31%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.556.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.217.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.775"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.724.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("E19").Value = "  -5.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.05%  "
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("E26").Value = "  -4.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "153.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("E32").Value = "  -8.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0756"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.65%  "
$ws.Range("E34").Value = "  -5.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.117"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("E36").Value = "  -8.44%  "
$ws.Range("E37").Value = "  -7.02%  "
$ws.Range("E38").Value = "  -5.49%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("B41").Value = "NEARProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.14%  "
$ws.Range("E42").Value = "  -10.83%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.835.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.16%  "
$ws.Range("E46").Value = "  -5.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "67.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.80%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "13.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.32%  "
